$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 88: correct the timestamp in column A ---
# (was 2024-06-26 14:24:24.000  ->  now 2024-06-26 07:00:00.000)
$ws.Range("A88").Value = 45469.2916666667

# --- New row 89 (latest day of data from the R script run) ---

# Column A: date/time value. Copy the style from A88 first so the new
# cell picks up the same "yyyy-mm-dd hh:mm:ss" date formatting (s="1"),
# then set the numeric serial value.
$ws.Range("A88").Copy()
$ws.Range("A89").PasteSpecial(-4122)
$ws.Range("A89").Value = 45470.6494212963

$ws.Range("B89").Value = 13500
$ws.Range("C89").Value = 3.15000009536743
$ws.Range("D89").Value = 2.97000002861023
$ws.Range("E89").Value = 2.97000002861023
$ws.Range("F89").Value = 3.22000002861023

# Column G: text representation of the close price (mirrors the pattern
# used by every other row). The leading apostrophe forces text so it is
# stored as a shared string instead of being re-interpreted as a number.
$ws.Range("G89").Value = "'3.22000002861023"
$ws.Range("G89").Style = $ws.Range("G88").Style

# Column H: ticker symbol, same shared string used throughout the sheet.
$ws.Range("H89").Value = "ESPE.MI"
